# Auto-generated edit script applying crypto price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.015.50'
$ws.Range("E2").Value = '  +1.37%  '
$ws.Range("D3").Value = '2.319.33'
$ws.Range("E3").Value = '  +1.82%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '117.44'
$ws.Range("E5").Value = '  +24.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '270.86'
$ws.Range("E6").Value = '  +1.65%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("E7").Value = '  +1.34%  '
$ws.Range("E8").Value = '  +0.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.629'
$ws.Range("E9").Value = '  +4.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '49.67'
$ws.Range("E10").Value = '  +11.93%  '
$ws.Range("E11").Value = '  +1.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.76'
$ws.Range("E12").Value = '  +13.88%  '
$ws.Range("E13").Value = '  +2.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.82'
$ws.Range("E14").Value = '  +4.80%  '
$ws.Range("D15").Value = '2.630.41'
$ws.Range("E15").Value = '  +0.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.868'
$ws.Range("E16").Value = '  +2.97%  '
$ws.Range("D17").Value = '2.319.91'
$ws.Range("E17").Value = '  +1.68%  '
$ws.Range("D18").Value = '43.932.45'
$ws.Range("E18").Value = '  +1.15%  '
$ws.Range("E19").Value = '  +4.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.66'
$ws.Range("E20").Value = '  +8.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.80'
$ws.Range("E21").Value = '  +1.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.55'
$ws.Range("E22").Value = '  +7.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.18'
$ws.Range("E23").Value = '  +0.75%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.92'
$ws.Range("E24").Value = '  +17.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.60'
$ws.Range("E25").Value = '  +7.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.60'
$ws.Range("E27").Value = '  +3.24%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '44.13'
$ws.Range("E28").Value = '  +13.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.42'
$ws.Range("E29").Value = '  -0.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.28'
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '178.30'
$ws.Range("E31").Value = '  +1.96%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.97'
$ws.Range("E32").Value = '  +0.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0941'
$ws.Range("E33").Value = '  +6.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.62'
$ws.Range("E34").Value = '  +5.87%  '
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.79'
$ws.Range("E35").Value = '  +9.11%  '
$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.126'
$ws.Range("E36").Value = '  +1.18%  '
$ws.Range("E37").Value = '  +4.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.98'
$ws.Range("E38").Value = '  +21.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0360'
$ws.Range("E39").Value = '  +1.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.250'
$ws.Range("E40").Value = '  +6.47%  '
$ws.Range("B41").Value = 'MultiversX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '75.04'
$ws.Range("E41").Value = '  +20.66%  '
$ws.Range("B42").Value = 'LidoDAOToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.41'
$ws.Range("E42").Value = '  +3.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.38'
$ws.Range("E43").Value = '  +13.17%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.08'
$ws.Range("E44").Value = '  +17.12%  '
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("E46").Value = '  +5.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.87'
$ws.Range("E47").Value = '  +1.14%  '
$ws.Range("E48").Value = '  -0.54%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.20'
$ws.Range("E49").Value = '  +3.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.24'
$ws.Range("E50").Value = '  +4.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.462'
$ws.Range("E51").Value = '  +9.17%  '
